# Locate the final inline picture (the last drawing in the document); its
# Range marks the end of the paragraph that also carries the "_GoBack"
# bookmark.  Collapse to the end of that shape and split the paragraph
# there so the bookmark stays with the picture's paragraph, and the new
# content becomes its own paragraphs inserted just before the bookmark.
$d = $word.ActiveDocument

$lastShape = $d.InlineShapes.Item($d.InlineShapes.Count)
$splitPoint = $lastShape.Range
$splitPoint.Collapse(0)   # wdCollapseEnd
$splitPoint.InsertParagraphAfter()

# The freshly-created (still empty) paragraph is now the last paragraph in
# the document; turn it into the new "Heading 2" section title.
$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Style = "Heading 2"
$headingPara.Range.InsertAfter("Using UART_TEST VisualDSP++ project")

# Add a new paragraph after the heading for the body copy, and give it the
# default "Normal" style (matching the surrounding body paragraphs).
$headingRange = $headingPara.Range
$headingRange.Collapse(0)   # wdCollapseEnd
$headingRange.InsertParagraphAfter()

$bodyPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bodyPara.Style = "Normal"

$bodyText = "This project is located at dsp_pid.git\src\vdsp_proj\test_uart and contains a VisualDSP project to output the letter " + [char]0x2018 + "a" + [char]0x2019 + " at 9600,8,N,1 baud rate.  The debugger will halt when the processor receives any UART byte."
$bodyPara.Range.InsertAfter($bodyText)
